# RF classify (read data header bug)
#
# The sheet used to carry two parallel copies of the table: a "raw" block in
# columns H:M (absolute video timestamps) plus a "derived" block in columns
# A:F that subtracted the video start time (row 9 / B9) from the raw block
# via formulas. That duplicated header reader caused the bug described in
# the commit message. The fix removes the now-unneeded raw H:M block and the
# "Video start time:" helper row, keeping only the clean A:F table with its
# values (no more formulas pointing at the deleted raw block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Freeze the A:F formulas (which read H:M and B$9) down to plain values
#    before the cells they depend on go away.
$ws.Range("A2:B7").Value = $ws.Range("A2:B7").Value()

# 2) Remove the "Video start time:" helper row (old row 9) completely; this
#    shifts everything below it up by one row (old row 11 -> row 10).
$ws.Rows("9").Delete()

# 3) Remove the duplicated raw-data block in columns H:M (header row copy +
#    the per-row raw timestamps/category columns).
$ws.Range("H1:M10").Delete()

# 4) Restore a plausible active selection on the cleaned-up sheet.
$ws.Range("Q11").Select() | Out-Null
